$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.803.92"
$ws.Range("E2").Value = "  +1.90%  "
$ws.Range("D3").Value = "3.150.75"
$ws.Range("E3").Value = "  +3.35%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.84"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.148.23"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("E9").Value = "  +5.30%  "
$ws.Range("E10").Value = "  +6.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.505"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +8.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +15.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.44"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +11.07%  "
$ws.Range("D15").Value = "3.662.15"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "64.846.84"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.21"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +8.23%  "
$ws.Range("D18").Value = "3.146.91"
$ws.Range("E18").Value = "  +3.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.09"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +8.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.737"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +10.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +10.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.86"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.18%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +5.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.91"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +12.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.91"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +10.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +11.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.62"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.69"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "486.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +12.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0867"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.10%  "
$ws.Range("E39").Value = "  +5.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("D41").Value = "3.119.14"
$ws.Range("E41").Value = "  +6.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.67"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.66%  "
$ws.Range("E43").Value = "  +6.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.294"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +14.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +17.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.62"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.95%  "
$ws.Range("D47").Value = "0.0₃0576"
$ws.Range("E47").Value = "  +13.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +12.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.40"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.09%  "
